$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.386.14"
$ws.Range("E2").Value = "  +4.29%  "

$ws.Range("D3").Value = "3.486.24"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Formula = "'585.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("E6").Value = "  +7.45%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.57%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Formula = "'0.127"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.57%  "

$ws.Range("D11").Formula = "'0.398"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.32%  "

$ws.Range("D12").Value = "4.081.79"
$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").Formula = "'29.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.92%  "

$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "3.483.76"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("E16").Value = "  +3.81%  "

$ws.Range("D17").Value = "63.381.03"

$ws.Range("D18").Formula = "'6.29"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Formula = "'14.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.37%  "

$ws.Range("E20").Value = "  +5.57%  "

$ws.Range("D21").Formula = "'392.65"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("E22").Value = "  +3.22%  "

$ws.Range("D23").Formula = "'75.23"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  +8.40%  "

$ws.Range("D26").Value = "3.630.51"
$ws.Range("E26").Value = "  +3.77%  "

$ws.Range("D27").Formula = "'0.184"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.94%  "

$ws.Range("E28").Value = "  +9.66%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Formula = "'8.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.39%  "

$ws.Range("E31").Value = "  +2.37%  "

$ws.Range("D32").Formula = "'1.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.76%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Formula = "'23.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.76%  "

$ws.Range("D35").Formula = "'32.44"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +28.01%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Formula = "'7.16"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.14%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Formula = "'5.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.58%  "

$ws.Range("D38").Formula = "'171.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("E39").Value = "  +9.77%  "

$ws.Range("D40").Value = "3.523.58"
$ws.Range("E40").Value = "  +3.69%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("E43").Value = "  +7.40%  "

$ws.Range("E44").Value = "  +3.92%  "

$ws.Range("D45").Formula = "'42.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("E46").Value = "  +9.86%  "

$ws.Range("D47").Value = "2.617.68"
$ws.Range("E47").Value = "  +7.10%  "

$ws.Range("D48").Formula = "'23.92"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.24%  "

$ws.Range("D49").Formula = "'2.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +18.44%  "

$ws.Range("E50").Value = "  +2.13%  "

$ws.Range("D51").Formula = "'0.0270"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.25%  "
